# Finished Week 13 logging
$wb = $excel.ActiveWorkbook

# Sheet "OFF" (first sheet) - row 3 (label "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 183
$wsOff.Range("C3").Value = 127
$wsOff.Range("D3").Value = 53
$wsOff.Range("E3").Value = 24

# Sheet "DEF" (second sheet) - row 3 (label "R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 238
$wsDef.Range("C3").Value = 183
$wsDef.Range("D3").Value = 35
$wsDef.Range("E3").Value = 17
$wsDef.Range("G3").Value = 2
